$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cell F1 "time_taken", copying the style from E1 (bold/border/center)
$ws.Cells.Item(1, 5).Copy() | Out-Null
$ws.Cells.Item(1, 6).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Cells.Item(1, 6).Value = "time_taken"

# Fill F2:F83 with the time_taken values (plain inline/shared strings, no special style)
$timeTaken = @(
    "2021-10-05 10:52:24.080406",
    "2021-10-05 10:52:24.080419",
    "2021-10-05 10:52:24.080422",
    "2021-10-05 10:52:24.080426",
    "2021-10-05 10:52:24.080430",
    "2021-10-05 10:52:24.080433",
    "2021-10-05 10:52:24.080436",
    "2021-10-05 10:52:24.080438",
    "2021-10-05 10:52:24.080441",
    "2021-10-05 10:52:24.080444",
    "2021-10-05 10:52:24.080447",
    "2021-10-05 10:52:24.080450",
    "2021-10-05 10:52:24.080453",
    "2021-10-05 10:52:24.080456",
    "2021-10-05 10:52:24.080458",
    "2021-10-05 10:52:24.080461",
    "2021-10-05 10:52:24.080464",
    "2021-10-05 10:52:24.080467",
    "2021-10-05 10:52:24.080470",
    "2021-10-05 10:52:24.080472",
    "2021-10-05 10:52:24.080475",
    "2021-10-05 10:52:24.080478",
    "2021-10-05 10:52:24.080481",
    "2021-10-05 10:52:24.080484",
    "2021-10-05 10:52:24.080487",
    "2021-10-05 10:52:24.080491",
    "2021-10-05 10:52:24.080494",
    "2021-10-05 10:52:24.080497",
    "2021-10-05 10:52:24.080499",
    "2021-10-05 10:52:24.080502",
    "2021-10-05 10:52:24.080505",
    "2021-10-05 10:52:24.080507",
    "2021-10-05 10:52:24.080511",
    "2021-10-05 10:52:24.080514",
    "2021-10-05 10:52:24.080516",
    "2021-10-05 10:52:24.080519",
    "2021-10-05 10:52:24.080522",
    "2021-10-05 10:52:24.080524",
    "2021-10-05 10:52:24.080527",
    "2021-10-05 10:52:24.080530",
    "2021-10-05 10:52:24.080533",
    "2021-10-05 10:52:24.080536",
    "2021-10-05 10:52:24.080539",
    "2021-10-05 10:52:24.080542",
    "2021-10-05 10:52:24.080544",
    "2021-10-05 10:52:24.080547",
    "2021-10-05 10:52:24.080550",
    "2021-10-05 10:52:24.080553",
    "2021-10-05 10:52:24.080555",
    "2021-10-05 10:52:24.080558",
    "2021-10-05 10:52:24.080561",
    "2021-10-05 10:52:24.080563",
    "2021-10-05 10:52:24.080566",
    "2021-10-05 10:52:24.080569",
    "2021-10-05 10:52:24.080572",
    "2021-10-05 10:52:24.080575",
    "2021-10-05 10:52:24.080577",
    "2021-10-05 10:52:24.080581",
    "2021-10-05 10:52:24.080583",
    "2021-10-05 10:52:24.080586",
    "2021-10-05 10:52:24.080589",
    "2021-10-05 10:52:24.080591",
    "2021-10-05 10:52:24.080594",
    "2021-10-05 10:52:24.080597",
    "2021-10-05 10:52:24.080601",
    "2021-10-05 10:52:24.080604",
    "2021-10-05 10:52:24.080607",
    "2021-10-05 10:52:24.080610",
    "2021-10-05 10:52:24.080612",
    "2021-10-05 10:52:24.080615",
    "2021-10-05 10:52:24.080618",
    "2021-10-05 10:52:24.080620",
    "2021-10-05 10:52:24.080623",
    "2021-10-05 10:52:24.080626",
    "2021-10-05 10:52:24.080629",
    "2021-10-05 10:52:24.080632",
    "2021-10-05 10:52:24.080636",
    "2021-10-05 10:52:24.080639",
    "2021-10-05 10:52:24.080642",
    "2021-10-05 10:52:24.080644",
    "2021-10-05 10:52:24.080647",
    "2021-10-05 10:52:24.080650"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

